# "Added uren voor groepswer 17-02-2016"
# Week 2 (column D) "gemeenschappelijk" (group-work) row 10 gains 4 extra
# hours: 4:15 -> 8:15, i.e. 0.177083333... -> 0.34375 (fraction of a day,
# [h]:mm formatted). The weekly totals in row 12-15 are driven by formulas
# referencing D10, so they recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 0.34375

# Leave the selection where the user finished editing.
[void]$ws.Range("E10").Select()
